$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "purpose" column (E2:E25) from "S.GISH" to the new value "fullRNASEQ"
$ws.Range("E2:E25").Value = "fullRNASEQ"

# Reflect the final selection left after the edit
$ws.Range("E24:E25").Select()
